$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L header "add"
$ws.Range("L1").Value = "add"

# Fill L2:L147 with sequential values starting at 4 (L2=4, L3=5, ... L147=149)
for ($row = 2; $row -le 147; $row++) {
    $ws.Cells.Item($row, 12).Value = $row + 2
}

# Update the active selection to match the post-edit state
$ws.Range("N19").Select()
